# "Big update. Login fixed of typo's"
#
# Fill in Sprint (F) / Status (G) = "Done" for the userstories that are
# now finished, and flip the "Update Functie" legend swatch (E28) from
# "in progress" (orange) to "done" (green). Finally move the on-screen
# selection to where the author left off (G17), with the view scrolled
# back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 4 - Ivanka Novoretzj userstory: Sprint 24, Status Done
$ws.Range("F4").Value = 24
$ws.Range("G4").Value = "Done"

# Row 5 - Ae Ri Hwang userstory: Sprint 27, Status Done
$ws.Range("F5").Value = 27
$ws.Range("G5").Value = "Done"

# Row 9 - Imogen Allen userstory: Sprint 2, Status Done
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "Done"

# Row 10 - Isabella de Ruyter userstory: Sprint 8, Status Done
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = "Done"

# Row 11 - Desiderius N. userstory: Sprint 9, Status Done
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = "Done"

# Row 17 - Cerys Osborne userstory: Sprint 15, Status Done
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = "Done"

# Legend: "Update Functie" swatch (E28) goes from Oranje (in progress) to
# Groen (done) - copy the fill straight from an already-green swatch so the
# exact same style gets reused.
$ws.Range("E28").Interior.Color = $ws.Range("E27").Interior.Color

# Leave the view where the author left it: scrolled back to the top of the
# sheet with G17 (the last edited cell) selected.
$ws.Range("G17").Select()
